$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Range("B16").Value = 383
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)  # xlPasteFormats
